$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V3")

# Row 5 - Room 2 1st move
$ws.Range("A5").Value = "Room 2 1st move"
$ws.Range("C5").Value = 1650
$ws.Range("B5").Value = 1553

# Row 6 - Room 2 Get key
$ws.Range("A6").Value = "Room 2 Get key"
$ws.Range("C6").Value = 1989
$ws.Range("B6").Value = 1899

# Row 8 - Room 3 1st move (added before row 7 so shared-string order matches source)
$ws.Range("A8").Value = "Room 3 1st move"
$ws.Range("C8").Value = 2833
$ws.Range("B8").Value = 2740

# Row 7 - Room 2 Enter door
$ws.Range("A7").Value = "Room 2 Enter door"
$ws.Range("C7").Value = 2116
$ws.Range("B7").Value = 2023

$wb.Application.Calculate()

$ws.Activate()
$ws.Range("A9").Select()
